{"js": "// MALS-1147 Renewal notice updates\n// Append \", Food and Fisheries\" right after \"Ministry of Agriculture\"\n// in the ministry heading paragraph near the bottom of the notice.\n\nconst body = context.document.body;\n\n// Locate the exact \"Ministry of Agriculture\" text (it appears exactly once,\n// in the Heading4-styled ministry-name paragraph at the foot of the notice).\nconst results = body.search(\"Ministry of Agriculture\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Ministry of Agriculture\" in the document body.');\n}\n\n// Insert the new text immediately after the matched text, extending the\n// paragraph to read \"Ministry of Agriculture, Food and Fisheries\".\nconst target = results.items[0];\ntarget.insertText(\", Food and Fisheries\", \"End\");\n\nawait context.sync();\n", "ps1": "# MALS-1147 Renewal notice updates\n# Append \", Food and Fisheries\" right after \"Ministry of Agriculture\"\n# in the ministry heading paragraph near the bottom of the notice.\n\n$d = $word.ActiveDocument\n\n# Locate the exact \"Ministry of Agriculture\" text (it appears exactly once,\n# in the Heading4-styled ministry-name paragraph at the foot of the notice).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Ministry of Agriculture\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\nif ($found -and $rng.Find.Found) {\n    # Collapse the found range to its end point, then insert the new text\n    # so the paragraph reads \"Ministry of Agriculture, Food and Fisheries\".\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\", Food and Fisheries\")\n}\n"}
